# Updated cryptos list on Wed Dec 20 21:42:10 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) snapshot columns for the coin
# table, and fixes rows 40/41 where the ranking order of NEARProtocol and
# LidoDAOToken swapped (their Coin name / Link / Price / Volume all move
# together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 40 and 41 swapped places in the ranking: row 40 is now LidoDAOToken,
# row 41 is now NEARProtocol.
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"

# Price (column D) and Volume(1h) (column E) refreshed values, keyed by row.
# These are stored as plain text in the workbook (e.g. "258.34", "43.695.89"
# which uses '.' thousands separators, and "  +3.08%  " padded with spaces),
# so every write below forces a Text number format first to stop the COM
# layer from coercing numeric-looking strings (like "258.34") into floating
# point numbers.
$updates = @{
    2  = @{ D = "43.695.89";  E = "  +3.08%  " }
    3  = @{ D = "2.188.06";   E = "  +0.29%  " }
    4  = @{ E = "  +0.01%  " }
    5  = @{ D = "258.34";     E = "  +1.74%  " }
    6  = @{ D = "81.06";      E = "  +9.89%  " }
    7  = @{ E = "  +1.65%  " }
    8  = @{ E = "  -0.01%  " }
    9  = @{ D = "0.592";      E = "  +1.17%  " }
    10 = @{ D = "42.84";      E = "  +6.04%  " }
    11 = @{ D = "0.0916";     E = "  +0.15%  " }
    12 = @{ D = "6.95";       E = "  +2.53%  " }
    13 = @{ E = "  +2.01%  " }
    14 = @{ D = "2.514.44";   E = "  +0.14%  " }
    15 = @{ D = "14.22";      E = "  +0.40%  " }
    16 = @{ D = "2.231.80";   E = "  +1.86%  " }
    17 = @{ D = "0.776";      E = "  +0.67%  " }
    18 = @{ D = "43.599.08";  E = "  +2.97%  " }
    19 = @{ D = "0.0000102";  E = "  +0.75%  " }
    20 = @{ D = "69.97";      E = "  -0.88%  " }
    21 = @{ D = "5.92";       E = "  +0.78%  " }
    22 = @{ D = "2.41";       E = "  +13.20%  " }
    23 = @{ D = "229.98";     E = "  +1.37%  " }
    24 = @{ D = "8.88";       E = "  -5.58%  " }
    25 = @{ E = "  +0.17%  " }
    26 = @{ D = "41.95";      E = "  +14.88%  " }
    27 = @{ D = "10.64";      E = "  +1.80%  " }
    28 = @{ E = "  +0.00%  " }
    29 = @{ D = "2.24";       E = "  +1.09%  " }
    30 = @{ D = "2.22";       E = "  +2.17%  " }
    31 = @{ D = "173.03";     E = "  +1.41%  " }
    32 = @{ D = "20.33";      E = "  +1.59%  " }
    33 = @{ D = "0.0871";     E = "  +8.37%  " }
    34 = @{ D = "5.27";       E = "  +3.02%  " }
    35 = @{ E = "  +5.14%  " }
    36 = @{ E = "  +1.44%  " }
    37 = @{ D = "4.45";       E = "  +4.97%  " }
    38 = @{ D = "0.0351";     E = "  +3.98%  " }
    39 = @{ D = "13.20";      E = "  +11.48%  " }
    40 = @{ D = "2.09";       E = "  +1.74%  " }
    41 = @{ D = "2.77";       E = "  +12.79%  " }
    42 = @{ D = "62.79";      E = "  +5.70%  " }
    43 = @{ D = "5.43";       E = "  +5.76%  " }
    44 = @{ E = "  +0.79%  " }
    45 = @{ D = "100.67";     E = "  -1.50%  " }
    46 = @{ D = "0.0983";     E = "  +1.34%  " }
    47 = @{ D = "8.20";       E = "  -0.51%  " }
    48 = @{ E = "  +4.25%  " }
    49 = @{ D = "1.54";       E = "  +25.86%  " }
    50 = @{ D = "1.10";       E = "  +1.37%  " }
    51 = @{ D = "0.436";      E = "  -6.49%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["E"]
    }
}
